$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 612.25
$ws.Range("I12").Value = 516.6667
$ws.Range("J12").Value = 899
$ws.Range("K12").Value = 516.6667
$ws.Range("L12").Value = 899
$ws.Range("M12").Value = -346.6667
$ws.Range("N12").Value = -1239
$ws.Range("H62").Value = 21856.285
$ws.Range("I62").Value = 22999
$ws.Range("K62").Value = 22999
$ws.Range("M62").Value = -22375
$ws.Range("H65").Value = 21856.285
$ws.Range("I65").Value = 22999
$ws.Range("K65").Value = 114995
$ws.Range("M65").Value = -111875
$ws.Range("H86").Value = 2388.889
$ws.Range("H89").Value = 2388.889
$ws.Range("H101").Value = 1396.4166
$ws.Range("I101").Value = 1727.75
$ws.Range("J101").Value = 733.75
$ws.Range("K101").Value = 5183.25
$ws.Range("L101").Value = 2201.25
$ws.Range("M101").Value = -3561.25
$ws.Range("N101").Value = -5445.25
$ws.Range("H113").Value = 21087.36
$ws.Range("I113").Value = 20982.54
$ws.Range("J113").Value = 21200.916
$ws.Range("K113").Value = 20982.54
$ws.Range("L113").Value = 21200.916
$ws.Range("M113").Value = -17728.54
$ws.Range("N113").Value = -27708.916
$ws.Range("H137").Value = 14703.737
$ws.Range("I137").Value = 5704.7896
$ws.Range("J137").Value = 23702.684
$ws.Range("K137").Value = 17114.3688
$ws.Range("L137").Value = 71108.052
$ws.Range("M137").Value = -14564.3688
$ws.Range("N137").Value = -76208.052
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1102549.9
$ws.Range("I32").Value = 1618310.1
$ws.Range("J32").Value = 71029.67999999999
$ws.Range("K32").Value = 1618310.1
$ws.Range("L32").Value = 71029.67999999999
$ws.Range("M32").Value = -1618023.1
$ws.Range("N32").Value = -71603.67999999999
$ws.Range("H45").Value = 3078.3447
$ws.Range("I45").Value = 3363.6365
$ws.Range("J45").Value = 2181.7144
$ws.Range("K45").Value = 3363.6365
$ws.Range("L45").Value = 2181.7144
$ws.Range("M45").Value = -2986.6365
$ws.Range("N45").Value = -2935.7144
$ws.Range("H74").Value = 13257.22
$ws.Range("I74").Value = 3422.1072
$ws.Range("J74").Value = 34440.54
$ws.Range("K74").Value = 3422.1072
$ws.Range("L74").Value = 34440.54
$ws.Range("M74").Value = -2548.1072
$ws.Range("N74").Value = -36188.54
$ws.Range("H77").Value = 13257.22
$ws.Range("I77").Value = 3422.1072
$ws.Range("J77").Value = 34440.54
$ws.Range("K77").Value = 17110.536
$ws.Range("L77").Value = 172202.7
$ws.Range("M77").Value = -12742.536
$ws.Range("N77").Value = -180938.7
$ws.Range("H102").Value = 24790.9
$ws.Range("J102").Value = 38416.5
$ws.Range("L102").Value = 38416.5
$ws.Range("N102").Value = -41660.5
$ws.Range("H122").Value = 6359.278
$ws.Range("I122").Value = 3922.5
$ws.Range("K122").Value = 11767.5
$ws.Range("M122").Value = -9317.5
$ws.Range("H133").Value = 57866.332
$ws.Range("J133").Value = 57866.332
$ws.Range("L133").Value = 57866.332
$ws.Range("N133").Value = -62926.332
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 1385.4667
$ws.Range("I11").Value = 203.55556
$ws.Range("J11").Value = 3158.3333
$ws.Range("K11").Value = 203.55556
$ws.Range("L11").Value = 3158.3333
$ws.Range("M11").Value = -63.55556000000001
$ws.Range("N11").Value = -3438.3333
$ws.Range("H20").Value = 15793.615
$ws.Range("I20").Value = 3659.1292
$ws.Range("J20").Value = 33706.43
$ws.Range("K20").Value = 3659.1292
$ws.Range("L20").Value = 33706.43
$ws.Range("M20").Value = -3412.1292
$ws.Range("N20").Value = -34200.43
$ws.Range("H94").Value = 5925.2856
$ws.Range("I94").Value = 2746.3333
$ws.Range("J94").Value = 24999
$ws.Range("K94").Value = 2746.3333
$ws.Range("L94").Value = 24999
$ws.Range("M94").Value = -2295.3333
$ws.Range("N94").Value = -25901
$ws.Range("H97").Value = 9739
$ws.Range("I97").Value = 9739
$ws.Range("K97").Value = 9739
$ws.Range("M97").Value = -8748
$ws.Range("H109").Value = 49995
$ws.Range("J109").Value = 49995
$ws.Range("L109").Value = 49995
$ws.Range("N109").Value = -52769
$ws.Range("H134").Value = 8577.159
$ws.Range("I134").Value = 2755.48
$ws.Range("K134").Value = 8266.440000000001
$ws.Range("M134").Value = -5731.440000000001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1297.8096
$ws.Range("I107").Value = 1146.6923
$ws.Range("K107").Value = 1146.6923
$ws.Range("M107").Value = 773.3077000000001
$ws.Range("H122").Value = 3619.92
$ws.Range("I122").Value = 3174.0435
$ws.Range("K122").Value = 9522.130500000001
$ws.Range("M122").Value = -7072.130500000001
$ws.Range("H132").Value = 5451.5293
$ws.Range("I132").Value = 1607.6774
$ws.Range("J132").Value = 11409.5
$ws.Range("K132").Value = 4823.0322
$ws.Range("L132").Value = 34228.5
$ws.Range("M132").Value = -2293.0322
$ws.Range("N132").Value = -39288.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1952.1904
$ws.Range("J132").Value = 2000
$ws.Range("L132").Value = 18000
$ws.Range("N132").Value = -23060
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 11000
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 11000
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 11000
$ws.Range("M15").Value = $null
$ws.Range("N15").Value = -11576
$ws.Range("H81").Value = 11000
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 11000
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 11000
$ws.Range("M81").Value = $null
$ws.Range("N81").Value = -12996
$ws.Range("H84").Value = 11000
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 11000
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 33000
$ws.Range("M84").Value = $null
$ws.Range("N84").Value = -42984
$ws.Range("H97").Value = 820.4375
$ws.Range("I97").Value = 780.9286
$ws.Range("K97").Value = 780.9286
$ws.Range("M97").Value = -284.9286
$ws.Range("H122").Value = 5446.625
$ws.Range("I122").Value = 3980
$ws.Range("J122").Value = 5935.5
$ws.Range("K122").Value = 11940
$ws.Range("L122").Value = 17806.5
$ws.Range("M122").Value = -9490
$ws.Range("N122").Value = -22706.5
$ws.Range("H132").Value = 4092.0876
$ws.Range("I132").Value = 2100.74
$ws.Range("K132").Value = 6302.219999999999
$ws.Range("M132").Value = -3772.219999999999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7378.2
$ws.Range("J7").Value = 8320.933999999999
$ws.Range("L7").Value = 8320.933999999999
$ws.Range("N7").Value = -8544.933999999999
$ws.Range("H122").Value = 7810.5713
$ws.Range("I122").Value = 5558.3335
$ws.Range("J122").Value = 9499.75
$ws.Range("K122").Value = 16675.0005
$ws.Range("L122").Value = 28499.25
$ws.Range("M122").Value = -14225.0005
$ws.Range("N122").Value = -33399.25
$ws.Range("H126").Value = 7378.2
$ws.Range("J126").Value = 8320.933999999999
$ws.Range("L126").Value = 24962.802
$ws.Range("N126").Value = -29902.802
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 53947.375
$ws.Range("J64").Value = 53947.375
$ws.Range("L64").Value = 53947.375
$ws.Range("N64").Value = -54443.375
$ws.Range("H67").Value = 53947.375
$ws.Range("J67").Value = 53947.375
$ws.Range("L67").Value = 53947.375
$ws.Range("N67").Value = -55663.375
$ws.Range("H96").Value = 2921.4443
$ws.Range("I96").Value = 2598.5
$ws.Range("K96").Value = 2598.5
$ws.Range("M96").Value = -1225.5
$ws.Range("H122").Value = 3201.8386
$ws.Range("I122").Value = 2306.5417
$ws.Range("J122").Value = 6271.4287
$ws.Range("K122").Value = 6919.625100000001
$ws.Range("L122").Value = 18814.2861
$ws.Range("M122").Value = -4469.625100000001
$ws.Range("N122").Value = -23714.2861
